# chore: update Sheets via scheduled runner
# Refreshes the computed Leve-profit columns (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> cols H-N)
# for a handful of rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with freshly pulled market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1596.3
$ws.Range("I40").Value = 1250.25
$ws.Range("J40").Value = 1827
$ws.Range("K40").Value = 1250.25
$ws.Range("L40").Value = 1827
$ws.Range("M40").Value = -1075.25
$ws.Range("N40").Value = -2177
$ws.Range("H64").Value = 45457456
$ws.Range("I64").Value = 166668910
$ws.Range("J64").Value = 3160
$ws.Range("K64").Value = 166668910
$ws.Range("L64").Value = 3160
$ws.Range("M64").Value = -166668662
$ws.Range("N64").Value = -3656
$ws.Range("H67").Value = 45457456
$ws.Range("I67").Value = 166668910
$ws.Range("J67").Value = 3160
$ws.Range("K67").Value = 166668910
$ws.Range("L67").Value = 3160
$ws.Range("M67").Value = -166668052
$ws.Range("N67").Value = -4876
$ws.Range("H76").Value = 6144.722
$ws.Range("I76").Value = 5176
$ws.Range("J76").Value = 6421.5
$ws.Range("K76").Value = 5176
$ws.Range("L76").Value = 6421.5
$ws.Range("M76").Value = -4861
$ws.Range("N76").Value = -7051.5
$ws.Range("H79").Value = 6144.722
$ws.Range("I79").Value = 5176
$ws.Range("J79").Value = 6421.5
$ws.Range("K79").Value = 5176
$ws.Range("L79").Value = 6421.5
$ws.Range("M79").Value = -4084
$ws.Range("N79").Value = -8605.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1335
$ws.Range("I63").Value = 1335
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1335
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -649
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1335
$ws.Range("I66").Value = 1335
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 6675
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -3243
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1710.8334
$ws.Range("I105").Value = 1335
$ws.Range("J105").Value = 2086.6667
$ws.Range("K105").Value = 1335
$ws.Range("L105").Value = 2086.6667
$ws.Range("M105").Value = 412
$ws.Range("N105").Value = -5580.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 22500
$ws.Range("J62").Value = 3750
$ws.Range("L62").Value = 3750
$ws.Range("N62").Value = -4998
$ws.Range("H65").Value = 22500
$ws.Range("J65").Value = 3750
$ws.Range("L65").Value = 18750
$ws.Range("N65").Value = -24990
$ws.Range("H132").Value = 7878.25
$ws.Range("I132").Value = 2100
$ws.Range("J132").Value = 9804.333000000001
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 29412.999
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -34472.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 808.97144
$ws.Range("I5").Value = 497.92
$ws.Range("J5").Value = 1586.6
$ws.Range("K5").Value = 1493.76
$ws.Range("L5").Value = 4759.799999999999
$ws.Range("M5").Value = -1381.76
$ws.Range("N5").Value = -4983.799999999999
$ws.Range("H122").Value = 1194.7
$ws.Range("I122").Value = 470.1111
$ws.Range("J122").Value = 1787.5454
$ws.Range("K122").Value = 4230.9999
$ws.Range("L122").Value = 16087.9086
$ws.Range("M122").Value = -1780.9999
$ws.Range("N122").Value = -20987.9086
$ws.Range("H130").Value = 2088.2856
$ws.Range("H131").Value = 551.4286
$ws.Range("I131").Value = 213.33333
$ws.Range("J131").Value = 1160
$ws.Range("K131").Value = 639.99999
$ws.Range("L131").Value = 3480
$ws.Range("M131").Value = 4400.00001
$ws.Range("N131").Value = -13560
$ws.Range("H135").Value = 808.97144
$ws.Range("I135").Value = 497.92
$ws.Range("J135").Value = 1586.6
$ws.Range("K135").Value = 4481.28
$ws.Range("L135").Value = 14279.4
$ws.Range("M135").Value = -1946.28
$ws.Range("N135").Value = -19349.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1053336
$ws.Range("I11").Value = 5000000
$ws.Range("J11").Value = 264003.2
$ws.Range("K11").Value = 5000000
$ws.Range("L11").Value = 264003.2
$ws.Range("M11").Value = -4999861
$ws.Range("N11").Value = -264281.2
$ws.Range("H42").Value = 41949
$ws.Range("J42").Value = 41949
$ws.Range("L42").Value = 41949
$ws.Range("N42").Value = -42919
$ws.Range("H70").Value = 4258.857
$ws.Range("I70").Value = 4078
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 4078
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -3808
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 4258.857
$ws.Range("I73").Value = 4078
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 4078
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -3142
$ws.Range("N73").Value = -6372
$ws.Range("H80").Value = 4875
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4875
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4875
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6871
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 4875
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4875
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 24375
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -34359
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H115").Value = 41949
$ws.Range("J115").Value = 41949
$ws.Range("L115").Value = 41949
$ws.Range("N115").Value = -44299
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4833923.5
$ws.Range("I10").Value = 14500000
$ws.Range("J10").Value = 885
$ws.Range("K10").Value = 14500000
$ws.Range("L10").Value = 885
$ws.Range("M10").Value = -14499860
$ws.Range("N10").Value = -1165
$ws.Range("H12").Value = 600
$ws.Range("I12").Value = 600
$ws.Range("K12").Value = 600
$ws.Range("M12").Value = -430
$ws.Range("H132").Value = 2423.195
$ws.Range("I132").Value = 2121.4075
$ws.Range("J132").Value = 3005.2144
$ws.Range("K132").Value = 6364.2225
$ws.Range("L132").Value = 9015.643199999999
$ws.Range("M132").Value = -3834.2225
$ws.Range("N132").Value = -14075.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 900
$ws.Range("I13").Value = 800
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 800
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -660
$ws.Range("N13").Value = -1280
